$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading percent values for rows 2-25, columns C,D,E,F,G,H,N
$data = @{
    "2" = @{ "C"=4.928768402600219; "D"=5.649667461805564; "E"=16.72339935239674; "F"=38.60592147596483; "G"=59.47479272158546; "H"=18.02939739762821; "N"=18.99769123883711 }
    "3" = @{ "C"=4.750715973229624; "D"=5.449517218026374; "E"=15.7413072787012; "F"=36.83191499190448; "G"=56.1154442933387; "H"=17.50004222806566; "N"=18.40031508502701 }
    "4" = @{ "C"=4.640201828005387; "D"=5.325647932214436; "E"=15.11449299435108; "F"=35.71611540401505; "G"=53.97312806182256; "H"=17.17499675511618; "N"=18.02485520896362 }
    "5" = @{ "C"=4.594943887929177; "D"=5.275008664708738; "E"=14.85333479962915; "F"=35.25546152999943; "G"=53.0810818156808; "H"=17.04274510246174; "N"=17.86990355188116 }
    "6" = @{ "C"=4.587417753661792; "D"=5.266592867242475; "E"=14.80963250433646; "F"=35.17863344159634; "G"=52.93184151972974; "H"=17.02080368120333; "N"=17.84406337566597 }
    "7" = @{ "C"=4.639592259107009; "D"=5.324965532006997; "E"=15.11099372557849; "F"=35.70992596150577; "G"=53.96117322644605; "H"=17.17321203841508; "N"=18.02277304766463 }
    "8" = @{ "C"=4.86766677092459; "D"=5.580906729155421; "E"=16.38985453637524; "F"=38.00015848500841; "G"=58.33356360182439; "H"=17.84700391936332; "N"=18.79364780656866 }
    "9" = @{ "C"=5.302440707935506; "D"=6.071699395578773; "E"=18.86970573823708; "F"=42.25359736752887; "G"=66.2409342593611; "H"=19.1593731761759; "N"=20.2273683202997 }
    "10" = @{ "C"=5.610646901156882; "D"=6.421499068240125; "E"=20.60686657495012; "F"=45.20348272566338; "G"=71.6079215133053; "H"=20.10724551707756; "N"=21.22223697909767 }
    "11" = @{ "C"=5.747787228625371; "D"=6.577581296934976; "E"=21.35638126211074; "F"=46.50257961559642; "G"=73.94827400440393; "H"=20.5328891772316; "N"=21.66018057919901 }
    "12" = @{ "C"=5.799233576771741; "D"=6.636197990370052; "E"=21.6344058659079; "F"=46.98803801519082; "G"=74.81966114465946; "H"=20.69312276344071; "N"=21.82377585681918 }
    "13" = @{ "C"=5.788175943633648; "D"=6.623596289287736; "E"=21.57478535990529; "F"=46.88377867963634; "G"=74.63265710415709; "H"=20.65865801399339; "N"=21.78864458690801 }
    "14" = @{ "C"=5.752029722615101; "D"=6.582413779872399; "E"=21.37937043541522; "F"=46.54265012011396; "G"=74.02026275025028; "H"=20.54609150817867; "N"=21.67368539489661 }
    "15" = @{ "C"=5.729824625547321; "D"=6.557123326602168; "E"=21.25891946079913; "F"=46.33284632113661; "G"=73.64321121104939; "H"=20.47701354564053; "N"=21.60297336126124 }
    "16" = @{ "C"=5.601618809704048; "D"=6.411233009933167; "E"=20.55706926943449; "F"=45.1176930226386; "G"=71.45291064196766; "H"=20.07930384196959; "N"=21.19330956972086 }
    "17" = @{ "C"=5.522149913736929; "D"=6.32091643456298; "E"=20.11611526749974; "F"=44.36101940554361; "G"=70.08308494553559; "H"=19.83378887045343; "N"=20.93814219015166 }
    "18" = @{ "C"=5.476155396883579; "D"=6.268684870150366; "E"=19.85865322088638; "F"=43.92178569123351; "G"=69.28569292124584; "H"=19.69206097939129; "N"=20.79000725568142 }
    "19" = @{ "C"=5.460534774780691; "D"=6.250953056740764; "E"=19.77082028624366; "F"=43.77238993259511; "G"=69.01408814294462; "H"=19.64399087572986; "N"=20.73962067985786 }
    "20" = @{ "C"=5.530639469913219; "D"=6.330560553955119; "E"=20.16345227902725; "F"=44.44198667598492; "G"=70.2298914315605; "H"=19.85997871186794; "N"=20.96544799483449 }
    "21" = @{ "C"=5.762660256602917; "D"=6.59452370368845; "E"=21.43692550303055; "F"=46.64302604969158; "G"=74.20054304427929; "H"=20.57918188555943; "N"=21.70751365554063 }
    "22" = @{ "C"=5.911451246215831; "D"=6.7641739590201; "E"=22.23543168887867; "F"=48.0436424859118; "G"=76.7089231171953; "H"=21.04363524480693; "N"=22.17935961385674 }
    "23" = @{ "C"=5.832312692421359; "D"=6.673905628839768; "E"=21.81232618168057; "F"=47.29966667446977; "G"=75.37816750397856; "H"=20.79630489046493; "N"=21.92877110911182 }
    "24" = @{ "C"=5.526802293663756; "D"=6.326201401030153; "E"=20.14206352909104; "F"=44.40539446720131; "G"=70.16355090943051; "H"=19.84814007350923; "N"=20.95310750188673 }
    "25" = @{ "C"=5.186557253758371; "D"=5.940557437116674; "E"=18.19349308199895; "F"=41.13180146318768; "G"=64.17740294065275; "H"=18.80647406668165; "N"=19.84905939529497 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
